$wb = $excel.ActiveWorkbook

# Update values on the "2010-Absolute" sheet
$ws = $wb.Worksheets.Item("2010-Absolute")
$ws.Range("C2").Value = 1000
$ws.Range("C3").Value = 600
$ws.Range("C5").Value = 0

# Move the active selection/tab from "2010-Relative" to "2010-Absolute"
$ws.Activate()
$ws.Range("C19").Select()

# Restore a plain selection on "2010-Relative" (still on cell D7, no longer the active tab)
$wsRel = $wb.Worksheets.Item("2010-Relative")
$wsRel.Range("D7").Select()

# Re-activate "2010-Absolute" so it remains the workbook's active/visible tab
$ws.Activate()
